$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4: change from text "2025-04-07 12:20:05" to the real datetime value,
# matching the numeric date formatting already used in B2/B3.
$ws.Range("B4").Value2 = 45754.51394675926
$ws.Range("B4").NumberFormat = "m/d/yy h:mm"

# Update the active cell selection to B6, as saved in the file.
$ws.Range("B6").Select()
